$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting from the neighboring
# header cell (G1) so it matches the bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"

# Era data updated: new Save column values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
